$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $idCell = $ws.Cells.Item($r, 1)
    $id = $idCell.Value2
    if ([string]::IsNullOrEmpty($id)) { continue }

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value2 = 45186
    }

    for ($col = 19; $col -le 25; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ([string]::IsNullOrEmpty($f)) { continue }
        if ($f -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $Matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $id + '")'
        }
    }
}
